$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells touched below to remain text (they are non-numeric,
# locale-formatted price strings in the source data) rather than being
# auto-converted to numbers by Excel when assigned.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.130.03'
$ws.Range("E2").Value = '  -2.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.89'
$ws.Range("E3").Value = '  -2.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.97'
$ws.Range("E5").Value = '  -1.77%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9989'
$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5126'
$ws.Range("E7").Value = '  +2.56%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3752'
$ws.Range("E8").Value = '  -1.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07136'
$ws.Range("E9").Value = '  -2.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8877'
$ws.Range("E10").Value = '  -2.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.65'
$ws.Range("E11").Value = '  -3.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.908.07'
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07543'
$ws.Range("E13").Value = '  -1.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.316'
$ws.Range("E14").Value = '  -2.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.16'
$ws.Range("E15").Value = '  -3.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9987'
$ws.Range("E16").Value = '  -0.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008475'
$ws.Range("E17").Value = '  -2.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.11'
$ws.Range("E18").Value = '  -3.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9988'
$ws.Range("E19").Value = '  +0.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.177.27'
$ws.Range("E20").Value = '  -2.71%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.047'
$ws.Range("E21").Value = '  -2.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.114.71'
$ws.Range("E22").Value = '  -0.25%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.56'
$ws.Range("E23").Value = '  -2.91%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.482'
$ws.Range("E24").Value = '  -1.82%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.83'
$ws.Range("E25").Value = '  -1.97%  '

$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.850'
$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.96'
$ws.Range("E27").Value = '  -2.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.100'
$ws.Range("E28").Value = '  -5.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.76'
$ws.Range("E29").Value = '  -1.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.722'
$ws.Range("E30").Value = '  -3.68%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.683'
$ws.Range("E31").Value = '  -2.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09018'
$ws.Range("E32").Value = '  +0.35%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05133'
$ws.Range("E33").Value = '  -2.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.084'
$ws.Range("E34").Value = '  -3.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7345'
$ws.Range("E36").Value = '  -6.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02053'
$ws.Range("E37").Value = '  -1.48%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.506'
$ws.Range("E38").Value = '  -5.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.057'
$ws.Range("E39").Value = '  -0.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.081'
$ws.Range("E40").Value = '  -0.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5348'
$ws.Range("E41").Value = '  -3.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.583'
$ws.Range("E42").Value = '  -3.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '117.00'
$ws.Range("E43").Value = '  +3.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.337'
$ws.Range("E44").Value = '  -2.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1475'
$ws.Range("E45").Value = '  -2.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4638'
$ws.Range("E46").Value = '  -3.89%  '

$ws.Range("E47").Value = '  +0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.08'
$ws.Range("E48").Value = '  -5.07%  '

$ws.Range("E49").Value = '  -4.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.42'
$ws.Range("E50").Value = '  -4.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.47'
$ws.Range("E51").Value = '  -1.56%  '
